$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 575 ("「融合／折衷／フュージョン」..."), shifting all
# subsequent rows up by one.
$ws.Rows.Item(575).Delete()
